$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update question text for each numbered row in column B (Sheet1).
# The underlying shared-strings table gets rebuilt on save, so we only
# need to set the final text per cell -- unused old strings are dropped
# and new/edited strings are appended automatically.
$ws.Range("B1").Value = " Do you feel comforatble around people"
$ws.Range("B2").Value = "Do you know how to captivate others"
$ws.Range("B3").Value = "Do you mind being center of attention"
$ws.Range("B4").Value = "Are you the life of the party?(Party animal)"
$ws.Range("B5").Value = "Do you keep yourself in background."
$ws.Range("B6").Value = "Do you wait for others to lead the way"
$ws.Range("B7").Value = "Are you quite around strangers?"
$ws.Range("B8").Value = "Do you bottle up your feelings "
$ws.Range("B9").Value = "Do you think you sympathise with others easily?"
$ws.Range("B10").Value = "Do you take out some of your time for others?"
$ws.Range("B11").Value = "Do you feel others emotions?"
$ws.Range("B12").Value = "Can you make others feel at ease"
$ws.Range("B13").Value = "Are you not really interested in other people's life?"
$ws.Range("B14").Value = "Do you insult others?"
$ws.Range("B15").Value = "Are you hard to get to know?"
$ws.Range("B16").Value = "Do you feel little concern for others."
$ws.Range("B17").Value = "Do you follow a regular schedule?"
$ws.Range("B18").Value = "Are you always prepared ?"
$ws.Range("B19").Value = "Are you exacting in your work?"
$ws.Range("B20").Value = "Do you pay attention to details?"
$ws.Range("B21").Value = "Do you waste your time?"
$ws.Range("B22").Value = "Often forget to put things back in their proper place."
$ws.Range("B23").Value = "Do you neglect your duties?"
$ws.Range("B24").Value = "Do things in half-way manner/"
$ws.Range("B25").Value = "Do you enjoy wild flights of fantasy?"
$ws.Range("B26").Value = "Enjoy thinking about things."
$ws.Range("B27").Value = "Believe in the importance of art"
$ws.Range("B28").Value = "Tend to vote for liberal political candidates."
$ws.Range("B29").Value = "Avoid philosophical discussions"
$ws.Range("B30").Value = " Do not like poetry?"
$ws.Range("B31").Value = "Rarely look for a deeper meaning in things."
$ws.Range("B32").Value = "Have difficulty understanding abstract ideas."
$ws.Range("B33").Value = "Seldom get mad?"
$ws.Range("B34").Value = "Are you filled with doubt about things?"
$ws.Range("B35").Value = "Have frequent mood swings."
$ws.Range("B36").Value = "Get stressed out easily"
$ws.Range("B37").Value = "Do you feel comfortable with yourself?"
$ws.Range("B38").Value = "Do you rarely lose your composure?"
$ws.Range("B39").Value = "Are you relaxed most of the time?"
$ws.Range("B40").Value = "Do you remain calm under pressure?"

# Restore the scrolled viewport position (was topLeftCell="A29", now "A27").
$ws.Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 27
$win.ScrollColumn = 1

# Keep the originally active cell selection (B37) intact.
$ws.Range("B37").Select()
